$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.161.89"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").Value = "3.526.46"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.90"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.23"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("D7").Value = "3.526.49"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.11"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.379"
$ws.Range("E12").Value = "  -1.88%  "
$ws.Range("D13").Value = "4.134.28"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.61"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.118"
$ws.Range("E15").Value = "  +1.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000179"
$ws.Range("E16").Value = "  -1.33%  "
$ws.Range("D17").Value = "3.517.00"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "64.247.89"
$ws.Range("E18").Value = "  -1.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.77"
$ws.Range("E19").Value = "  -3.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.04"
$ws.Range("E20").Value = "  -2.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.62"
$ws.Range("E21").Value = "  -1.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "385.22"
$ws.Range("E22").Value = "  -1.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.575"
$ws.Range("E23").Value = "  -0.98%  "
$ws.Range("D24").Value = "3.670.45"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.06"
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +3.11%  "
$ws.Range("E28").Value = "  -2.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.48"
$ws.Range("E29").Value = "  -3.48%  "
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("E31").Value = "  +0.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.23"
$ws.Range("E32").Value = "  -2.26%  "
$ws.Range("D33").Value = "3.537.05"
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.59"
$ws.Range("E35").Value = "  -2.05%  "
$ws.Range("E36").Value = "  +1.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.35"
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("E38").Value = "  -0.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.91"
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "160.97"
$ws.Range("E40").Value = "  -4.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0788"
$ws.Range("E41").Value = "  -2.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.815"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.54"
$ws.Range("E43").Value = "  +2.65%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.83"
$ws.Range("E45").Value = "  -2.43%  "
$ws.Range("E46").Value = "  -4.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.41"
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.61"
$ws.Range("E48").Value = "  -2.67%  "
$ws.Range("D49").Value = "2.474.02"
$ws.Range("E49").Value = "  +1.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.81"
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.907"
$ws.Range("E51").Value = "  -0.22%  "
